$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 3266.6667
$ws.Range("I51").Value = 3089.9
$ws.Range("J51").Value = 3427.3635
$ws.Range("K51").Value = 3089.9
$ws.Range("L51").Value = 3427.3635
$ws.Range("M51").Value = -2605.9
$ws.Range("N51").Value = -4395.363499999999
$ws.Range("H107").Value = 930.7273
$ws.Range("I107").Value = 523.8
$ws.Range("J107").Value = 5000
$ws.Range("K107").Value = 523.8
$ws.Range("L107").Value = 5000
$ws.Range("M107").Value = 1396.2
$ws.Range("N107").Value = -8840
$ws.Range("H138").Value = 1985
$ws.Range("J138").Value = 2046.5116
$ws.Range("L138").Value = 6139.5348
$ws.Range("N138").Value = -16419.5348

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 45455428
$ws.Range("I2").Value = 111111704
$ws.Range("J2").Value = 1086.7693
$ws.Range("K2").Value = 111111704
$ws.Range("L2").Value = 1086.7693
$ws.Range("M2").Value = -111111591
$ws.Range("N2").Value = -1312.7693
$ws.Range("H45").Value = 18519700
$ws.Range("I45").Value = 25642160
$ws.Range("J45").Value = 1302.8
$ws.Range("K45").Value = 25642160
$ws.Range("L45").Value = 1302.8
$ws.Range("M45").Value = -25641783
$ws.Range("N45").Value = -2056.8
$ws.Range("H61").Value = 1611.625
$ws.Range("I61").Value = 1155.8096
$ws.Range("J61").Value = 2481.818
$ws.Range("K61").Value = 1155.8096
$ws.Range("L61").Value = 2481.818
$ws.Range("M61").Value = -943.8096
$ws.Range("N61").Value = -2905.818
$ws.Range("H74").Value = 985.2353000000001
$ws.Range("I74").Value = 987.4375
$ws.Range("K74").Value = 987.4375
$ws.Range("M74").Value = -113.4375
$ws.Range("H77").Value = 985.2353000000001
$ws.Range("I77").Value = 987.4375
$ws.Range("K77").Value = 4937.1875
$ws.Range("M77").Value = -569.1875
$ws.Range("H102").Value = 1732.8572
$ws.Range("J102").Value = 1547.5
$ws.Range("L102").Value = 1547.5
$ws.Range("N102").Value = -4791.5
$ws.Range("H116").Value = 45455428
$ws.Range("I116").Value = 111111704
$ws.Range("J116").Value = 1086.7693
$ws.Range("K116").Value = 111111704
$ws.Range("L116").Value = 1086.7693
$ws.Range("M116").Value = -111109410
$ws.Range("N116").Value = -5674.7693
$ws.Range("H122").Value = 1182.9474
$ws.Range("I122").Value = 1070.1786
$ws.Range("J122").Value = 1498.7
$ws.Range("K122").Value = 3210.5358
$ws.Range("L122").Value = 4496.1
$ws.Range("M122").Value = -760.5357999999997
$ws.Range("N122").Value = -9396.1
$ws.Range("H132").Value = 5964.2
$ws.Range("I132").Value = 6844.5454
$ws.Range("K132").Value = 20533.6362
$ws.Range("M132").Value = -18003.6362
$ws.Range("H136").Value = 1611.625
$ws.Range("I136").Value = 1155.8096
$ws.Range("J136").Value = 2481.818
$ws.Range("K136").Value = 3467.4288
$ws.Range("L136").Value = 7445.454000000001
$ws.Range("M136").Value = -917.4288000000001
$ws.Range("N136").Value = -12545.454

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 45455428
$ws.Range("I3").Value = 111111704
$ws.Range("J3").Value = 1086.7693
$ws.Range("K3").Value = 111111704
$ws.Range("L3").Value = 1086.7693
$ws.Range("M3").Value = -111111590
$ws.Range("N3").Value = -1314.7693

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3791640.2
$ws.Range("I31").Value = 2851
$ws.Range("K31").Value = 2851
$ws.Range("M31").Value = -2556
$ws.Range("H34").Value = 3791640.2
$ws.Range("I34").Value = 2851
$ws.Range("K34").Value = 2851
$ws.Range("M34").Value = -2649
$ws.Range("H99").Value = 2635.125
$ws.Range("I99").Value = 2494.3333
$ws.Range("J99").Value = 2816.1428
$ws.Range("K99").Value = 2494.3333
$ws.Range("L99").Value = 2816.1428
$ws.Range("M99").Value = -996.3332999999998
$ws.Range("N99").Value = -5812.1428
$ws.Range("H126").Value = 2635.125
$ws.Range("I126").Value = 2494.3333
$ws.Range("J126").Value = 2816.1428
$ws.Range("K126").Value = 7482.999899999999
$ws.Range("L126").Value = 8448.428400000001
$ws.Range("M126").Value = -5012.999899999999
$ws.Range("N126").Value = -13388.4284
$ws.Range("H132").Value = 2978630.2
$ws.Range("I132").Value = 1806.7916
$ws.Range("J132").Value = 6947728.5
$ws.Range("K132").Value = 5420.3748
$ws.Range("L132").Value = 20843185.5
$ws.Range("M132").Value = -2890.3748
$ws.Range("N132").Value = -20848245.5
$ws.Range("H134").Value = 938.43475
$ws.Range("I134").Value = 842.94446
$ws.Range("J134").Value = 1282.2
$ws.Range("K134").Value = 2528.83338
$ws.Range("L134").Value = 3846.6
$ws.Range("M134").Value = 6.166619999999966
$ws.Range("N134").Value = -8916.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H46").Value = 563.5714
$ws.Range("J46").Value = 765.55554
$ws.Range("L46").Value = 2296.66662
$ws.Range("N46").Value = -2478.66662
$ws.Range("H96").Value = 11429.944
$ws.Range("I96").Value = 0
$ws.Range("J96").Value = 11429.944
$ws.Range("K96").Value = 0
$ws.Range("L96").Value = 34289.83199999999
$ws.Range("M96").ClearContents()
$ws.Range("N96").Value = -38407.83199999999
$ws.Range("H131").Value = 768.15
$ws.Range("J131").Value = 799.29346
$ws.Range("L131").Value = 2397.88038
$ws.Range("N131").Value = -12477.88038
$ws.Range("H132").Value = 1540.25
$ws.Range("I132").Value = 946.38464
$ws.Range("J132").Value = 2054.9333
$ws.Range("K132").Value = 8517.46176
$ws.Range("L132").Value = 18494.3997
$ws.Range("M132").Value = -5987.46176
$ws.Range("N132").Value = -23554.3997

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1523.6666
$ws.Range("I102").Value = 1100
$ws.Range("J102").Value = 1576.625
$ws.Range("K102").Value = 1100
$ws.Range("L102").Value = 1576.625
$ws.Range("M102").Value = 522
$ws.Range("N102").Value = -4820.625

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 8930522
$ws.Range("I68").Value = 1216.8334
$ws.Range("J68").Value = 15627500
$ws.Range("K68").Value = 1216.8334
$ws.Range("L68").Value = 15627500
$ws.Range("M68").Value = -467.8334
$ws.Range("N68").Value = -15628998
$ws.Range("H71").Value = 8930522
$ws.Range("I71").Value = 1216.8334
$ws.Range("J71").Value = 15627500
$ws.Range("K71").Value = 6084.166999999999
$ws.Range("L71").Value = 78137500
$ws.Range("M71").Value = -2340.166999999999
$ws.Range("N71").Value = -78144988
$ws.Range("H132").Value = 17500.084
$ws.Range("I132").Value = 37220.8
$ws.Range("K132").Value = 111662.4
$ws.Range("M132").Value = -109132.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 498.4375
$ws.Range("I107").Value = 384
$ws.Range("J107").Value = 612.875
$ws.Range("K107").Value = 1152
$ws.Range("L107").Value = 1838.625
$ws.Range("M107").Value = 768
$ws.Range("N107").Value = -5678.625
$ws.Range("H136").Value = 3504.1538
$ws.Range("I136").Value = 3387.8333
$ws.Range("K136").Value = 10163.4999
$ws.Range("M136").Value = -7613.499899999999
